$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("person")
$ws1.Range("A32:F35").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$v = $excel.ActiveWindow.ScrollRow
Write-Host "ScrollRow after set: $v"
